$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Funddatum (column D) for rows 30 and 31: 01.10.2020 -> 04.10.2020 (serial 44113 -> 44116)
$ws.Range("D30").Value = 44116
$ws.Range("D31").Value = 44116

# Add a new literature entry in row 32
# (cell entry order matches the shared-string insertion order of the target file)
$ws.Range("H32").Value = "https://books.google.de/books?hl=en&lr=&id=oZp4DwAAQBAJ&oi=fnd&pg=PR5&ots=CLnsbrW1xM&sig=t7eifDgNO5yFjCgxKXdoV3bC8LQ&redir_esc=y#v=onepage&q&f=false"
$ws.Range("A32").Value = "Softwarewartung: Grundlagen, Management und Wartungstechniken"
$ws.Range("B32").Value = "Christoph Bommer and Markus Spindler and Volkert Barr"
$ws.Range("G32").Value = "kaufen"
$ws.Range("E32").Value = "Beschreibt Grundlagen zur korrekten Wartung von Software (Es werden leider keine echten Methoden aufgezeigt)"
$ws.Range("C32").Value = 2008

# D32 needs the same "date" number format as D30/D31 (a new cell otherwise
# inherits the column's plain-number style) - copy the format from D31 first.
$ws.Range("D31").Copy()
$ws.Range("D32").PasteSpecial(-4122)
$ws.Range("D32").Value = 44116

$ws.Range("F32").Value = 3

# Match row 32 formatting to similarly wrapped rows (row height + style already applied via E-column wrap style)
$ws.Rows.Item(32).RowHeight = 43.2

# Update view state: scroll/freeze pane position and active selection
$ws.Application.ActiveWindow.ScrollRow = 18
$ws.Range("C31").Select()
